# Remove two account rows from the "Export" sheet:
#   - 004206790 / EMMANUELLE / 158000   (originally Excel row 2)
#   - 004556853 / MARCEL     / 2077.58  (originally Excel row 6)
#
# Delete from the bottom up so the earlier deletion doesn't shift the
# row index of the one still to be removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").EntireRow.Delete()
$ws.Range("A2").EntireRow.Delete()
